# Updates cryptos list (prices/volumes) and swaps the Mantle/Kaspa rows.
# Price-column (D) values are stored as TEXT in the source data (several
# contain two '.' separators, e.g. "69.680.66", which can't be numbers).
# Force text storage via NumberFormat "@" before assigning so Excel
# doesn't silently reinterpret plain-numeric-looking strings (e.g.
# "604.22") as floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "69.680.66"

# Row 3 - Ethereum
Set-TextValue "D3" "3.886.69"
$ws.Range("E3").Value = "  +1.31%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
Set-TextValue "D5" "604.22"
$ws.Range("E5").Value = "  +0.91%  "

# Row 6 - Solana
Set-TextValue "D6" "170.13"
$ws.Range("E6").Value = "  +4.93%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.885.89"
$ws.Range("E7").Value = "  +1.32%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.11%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +1.40%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.66%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +1.07%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +2.05%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +4.93%  "

# Row 14 - Avalanche
$ws.Range("E14").Value = "  +3.94%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.543.81"
$ws.Range("E15").Value = "  +1.47%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.880.54"
$ws.Range("E16").Value = "  +0.68%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "69.661.45"
$ws.Range("E17").Value = "  +1.56%  "

# Row 18 - Chainlink
Set-TextValue "D18" "18.67"
$ws.Range("E18").Value = "  +9.38%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +1.77%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  -0.69%  "

# Row 21 - Uniswap
Set-TextValue "D21" "11.17"
$ws.Range("E21").Value = "  -0.66%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "489.34"
$ws.Range("E22").Value = "  +1.28%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.748"
$ws.Range("E23").Value = "  +4.73%  "

# Row 24 - PEPE
$ws.Range("E24").Value = "  +3.31%  "

# Row 25 - Litecoin
Set-TextValue "D25" "85.28"
$ws.Range("E25").Value = "  +1.76%  "

# Row 26 - Fetch.AI
$ws.Range("E26").Value = "  +4.00%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "12.37"
$ws.Range("E27").Value = "  +2.52%  "

# Row 28 - RenderToken
Set-TextValue "D28" "10.13"
$ws.Range("E28").Value = "  +2.41%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.25%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +1.25%  "

# Row 31 - WrappedeETH
Set-TextValue "D31" "4.037.99"
$ws.Range("E31").Value = "  +1.27%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  +2.60%  "

# Row 33 - NEARProtocol
Set-TextValue "D33" "7.85"
$ws.Range("E33").Value = "  +0.59%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "31.87"
$ws.Range("E34").Value = "  -0.19%  "

# Row 35 - RenzoRestakedETH
Set-TextValue "D35" "3.854.91"
$ws.Range("E35").Value = "  +1.89%  "

# Row 36 - Hedera
$ws.Range("E36").Value = "  +0.68%  "

# Row 37 - Filecoin
Set-TextValue "D37" "6.12"
$ws.Range("E37").Value = "  +4.75%  "

# Row 38 / 39 - swap Mantle <-> Kaspa (Kaspa now ranked above Mantle)
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D38" "0.142"
$ws.Range("E38").Value = "  +2.34%  "

$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D39" "1.03"
$ws.Range("E39").Value = "  +0.74%  "

# Row 40 - dogwifhat
Set-TextValue "D40" "3.36"
$ws.Range("E40").Value = "  +14.12%  "

# Row 41 - FirstDigitalUSD
$ws.Range("E41").Value = "  +0.07%  "

# Row 42 - TheGraph
$ws.Range("E42").Value = "  +3.72%  "

# Row 43 - Stacks
$ws.Range("E43").Value = "  +5.93%  "

# Row 44 - Bittensor
Set-TextValue "D44" "437.81"
$ws.Range("E44").Value = "  +2.17%  "

# Row 45 - OKB
Set-TextValue "D45" "48.12"
$ws.Range("E45").Value = "  -0.77%  "

# Row 46 - Cosmos
$ws.Range("E46").Value = "  +4.04%  "

# Row 48 - FLOKI
Set-TextValue "D48" "0.000276"
$ws.Range("E48").Value = "  +22.95%  "

# Row 49 - VeChain
Set-TextValue "D49" "0.0367"
$ws.Range("E49").Value = "  +3.57%  "

# Row 51 - Arweave
Set-TextValue "D51" "40.52"
$ws.Range("E51").Value = "  +5.10%  "
